$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the new "cost" (column L) values for rows 2-60, formatted as integer ("0")
$ws.Range("L2:L60").NumberFormat = "0"

$ws.Range("L2").Value = 9
$ws.Range("L3").Value = 13
$ws.Range("L4").Value = 19
$ws.Range("L5").Value = 28
$ws.Range("L6").Value = 40
$ws.Range("L7").Value = 59
$ws.Range("L8").Value = 87
$ws.Range("L9").Value = 127
$ws.Range("L10").Value = 186
$ws.Range("L11").Value = 273
$ws.Range("L12").Value = 399
$ws.Range("L13").Value = 585
$ws.Range("L14").Value = 856
$ws.Range("L15").Value = 1254
$ws.Range("L16").Value = 1837
$ws.Range("L17").Value = 2691
$ws.Range("L18").Value = 3941
$ws.Range("L19").Value = 5772
$ws.Range("L20").Value = 8454
$ws.Range("L21").Value = 12382
$ws.Range("L22").Value = 18135
$ws.Range("L23").Value = 26561
$ws.Range("L24").Value = 38902
$ws.Range("L25").Value = 56977
$ws.Range("L26").Value = 83451
$ws.Range("L27").Value = 122226
$ws.Range("L28").Value = 179017
$ws.Range("L29").Value = 262195
$ws.Range("L30").Value = 384022
$ws.Range("L31").Value = 562454
$ws.Range("L32").Value = 823793
$ws.Range("L33").Value = 1206561
$ws.Range("L34").Value = 1767179
$ws.Range("L35").Value = 2588282
$ws.Range("L36").Value = 3790904
$ws.Range("L37").Value = 5552314
$ws.Range("L38").Value = 8132146
$ws.Range("L39").Value = 11910675
$ws.Range("L40").Value = 17444863
$ws.Range("L41").Value = 25550461
$ws.Range("L42").Value = 37422253
$ws.Range("L43").Value = 54810166
$ws.Range("L44").Value = 80277216
$ws.Range("L45").Value = 117577302
$ws.Range("L46").Value = 172208537
$ws.Range("L47").Value = 252223684
$ws.Range("L48").Value = 369417149
$ws.Range("L49").Value = 541063502
$ws.Range("L50").Value = 792463789
$ws.Range("L51").Value = 1160674956
$ws.Range("L52").Value = 1699972128
$ws.Range("L53").Value = 2000000000
$ws.Range("L54").Value = 3646734750
$ws.Range("L55").Value = 5341157231
$ws.Range("L56").Value = 7822877869
$ws.Range("L57").Value = 11457707664
$ws.Range("L58").Value = 16781428411
$ws.Range("L59").Value = 24578768089
$ws.Range("L60").Value = 35999071473

# Restore selection to L2:L60 (matches author's final selection state)
$ws.Range("L2:L60").Select()
